$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 195; existing rows 195..235 shift down to 196..236
$ws.Rows.Item(195).Insert()

# Populate the newly inserted row 195 with the new record's data.
# Columns A, B, C, E, F, G, H, I, O, R carry the same constant values used
# throughout this data block, so copy them from the row now just below (196).
$ws.Cells.Item(195, 1).Value2 = $ws.Cells.Item(196, 1).Value2   # A - Mercado ID
$ws.Cells.Item(195, 2).Value2 = $ws.Cells.Item(196, 2).Value2   # B - Mercado
$ws.Cells.Item(195, 3).Value2 = $ws.Cells.Item(196, 3).Value2   # C - Region
$ws.Cells.Item(195, 5).Value2 = $ws.Cells.Item(196, 5).Value2   # E - Codreg
$ws.Cells.Item(195, 6).Value2 = $ws.Cells.Item(196, 6).Value2   # F - Categoria ID
$ws.Cells.Item(195, 7).Value2 = $ws.Cells.Item(196, 7).Value2   # G - Categoria
$ws.Cells.Item(195, 8).Value2 = $ws.Cells.Item(196, 8).Value2   # H - Variedad
$ws.Cells.Item(195, 9).Value2 = $ws.Cells.Item(196, 9).Value2   # I - Calidad
$ws.Cells.Item(195, 15).Value2 = $ws.Cells.Item(196, 15).Value2 # O - Origen
$ws.Cells.Item(195, 18).Value2 = $ws.Cells.Item(196, 18).Value2 # R - Clasificacion

# New record-specific values
$ws.Cells.Item(195, 4).Value2 = 44694        # D - Fecha
$ws.Cells.Item(195, 10).Value2 = 120         # J - Volumen
$ws.Cells.Item(195, 11).Value2 = 550         # K - Precio minimo
$ws.Cells.Item(195, 12).Value2 = 600         # L - Precio maximo
$ws.Cells.Item(195, 13).Value2 = 575         # M - Precio promedio ponderado
$ws.Cells.Item(195, 14).Value2 = "$/atado 0,5 a 1 kilo"   # N - Unidad de comercializacion
$ws.Cells.Item(195, 16).Value2 = 575         # P - Precio $/Kg
$ws.Cells.Item(195, 17).Value2 = 1           # Q - Kg o Unidades
